$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old J and K columns entirely (no longer used)
$ws.Range("J1:K5").Clear()

# --- Row 2 (Time horizon) numeric ranges updated ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 25

# --- Row 3 (Desired growth) numeric ranges updated ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 11

# --- Row 4 (Fluctuations) numeric ranges updated ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 100

# --- Row 5 (Worst case) last range updated ---
$ws.Range("F5").Value = 100

# --- New H/I columns, entered in the order the strings were first authored ---
$ws.Range("H1").Value = "Metric"
$ws.Range("H2").Value = "Window of Data"
$ws.Range("H4").Value = "STD Deviation"
$ws.Range("H5").Value = "Max Dropdown"
$ws.Range("I5").Value = "Filter"
$ws.Range("I1").Value = "Function"
$ws.Range("I2").Value = "Average Calculation Range"
$ws.Range("H3").Value = "CARG"
$ws.Range("I3").Value = "Y-Axis"
$ws.Range("I4").Value = "X-Axis"

# Column widths for new columns H and I (best-fit sized to their contents)
$ws.Columns("H").ColumnWidth = 13
$ws.Columns("I").ColumnWidth = 22.166666666666668

# Update selection to I3, matching the author's final cursor position
$ws.Range("I3").Select()

$wb.Save()
